# Change cell B11 from the shared string "R40" to the text string "1".
#
# A plain  $ws.Range("B11").Value = "1"  would be auto-coerced to the
# *number* 1 (Excel's normal type-inference for a numeric-looking value),
# which is not what the target state needs - B11 must keep its text type
# ("t=s") and its existing style/number-format (General).
#
# To force a genuine text value without disturbing B11's own formatting,
# stage the text on a scratch cell via TEXT() (whose natural result type
# is already a string, not a number), copy just that value onto B11, then
# remove the scratch cell again.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Formula = "=TEXT(1,""0"")"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
